# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# New values were recalculated (std/mean regen) and written into the
# existing "K" column (G) for rows 2-39.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 7
    3  = 0
    4  = 1
    5  = 1
    6  = 4
    7  = 3
    8  = 3
    9  = 3
    10 = 2
    11 = 5
    12 = 1
    13 = 1
    14 = 7
    15 = 2
    16 = 1
    17 = 6
    18 = 7
    19 = 6
    20 = 5
    21 = 4
    22 = 0
    23 = 1
    24 = 4
    25 = 6
    26 = 4
    27 = 3
    28 = 2
    29 = 3
    30 = 5
    31 = 4
    32 = 3
    33 = 3
    34 = 5
    35 = 3
    36 = 6
    37 = 2
    38 = 1
    39 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
